$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.126.45"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "1.796.79"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.88"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.45"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0717"
$ws.Range("E10").Value = "  +4.50%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "2.055.41"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "1.799.83"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.74"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "34.152.47"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.14"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.65"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.96"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.54"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0521"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.51"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "1.411.58"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.644"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.945"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.33"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  +4.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.96"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0498"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.04"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "1.954.42"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.08"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +0.79%  "
